$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mapping of row -> (DAMSLTag, DialogAct) updates (columns I and J)
$updates = @{
    5   = @("sv", "Statement-opinion")
    12  = @("sd", "Statement-non-opinion")
    16  = @("sv", "Statement-opinion")
    18  = @("sv", "Statement-opinion")
    19  = @("sd", "Statement-non-opinion")
    34  = @("sv", "Statement-opinion")
    35  = @("sv", "Statement-opinion")
    70  = @("sv", "Statement-opinion")
    74  = @("ba", "Appreciation")
    84  = @("sv", "Statement-opinion")
    89  = @("sd", "Statement-non-opinion")
    93  = @("%", "Uninterpretable")
    99  = @("sv", "Statement-opinion")
    101 = @("sd", "Statement-non-opinion")
    104 = @("sv", "Statement-opinion")
    106 = @("ba", "Appreciation")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("I$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
}
